$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the hidden "_GoBack" bookmark from the end of the image paragraph
#    (paragraph 3) into the now-empty paragraph 2, and drop that paragraph's
#    Heading2 style so it becomes a plain empty paragraph.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$p2.Range.Style = "Normal"
$d.Bookmarks.Add("_GoBack", $p2.Range) | Out-Null

# ---------------------------------------------------------------------------
# 2) Re-stamp the picture's drawing anchor/edit ids (Word regenerates these
#    hex ids whenever the drawing is touched/re-saved).  Rebuild paragraph 3
#    (the image paragraph) verbatim, with the new ids and without the
#    bookmark that used to trail the run (it now lives in paragraph 2).
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(3).Range
$p3.InsertXML('<w:p w:rsidR="0082648A" w:rsidRDefault="000928DE" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="7F1251C6" wp14:editId="7CD822CB"><wp:extent cx="2433955" cy="3657600"/><wp:effectExtent l="0" t="0" r="4445" b="0"/><wp:docPr id="1" name="Picture 1"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 1"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId5"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="2433955" cy="3657600"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>') | Out-Null

# ---------------------------------------------------------------------------
# 3) Materialize the built-in "Caption" style into the style sheet (it is
#    otherwise latent/unused - unhideWhenUsed - so apply it momentarily to a
#    throw-away range and then restore that range's original style).
# ---------------------------------------------------------------------------
$target = $d.Paragraphs(1)
$savedStyle = $target.Range.Style
$target.Range.Style = $d.Styles(-35)
$target.Range.Style = $savedStyle

$caption = $d.Styles("Caption")
$caption.NameLocal = "caption"
$caption.NextParagraphStyle = "Normal"
$caption.Priority = 35
$caption.UnhideWhenUsed = $true
$caption.Font.Bold = $true
$caption.Font.BoldBi = $true
$caption.Font.Size = 9
$caption.Font.SizeBi = 9
$caption.Font.TextColor.ObjectThemeColor = 4
$caption.ParagraphFormat.SpaceAfter = 10

Write-Output "done"
